$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A2: "12" -> "2025-04-29" (keep as text, not a date serial)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-04-29"
$ws.Range("A2").Style = "Normal"

# C2: "12" -> "123" (keep as text, not a number)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "123"
$ws.Range("C2").Style = "Normal"

# F2: "12" -> "21" (keep as text, not a number)
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "21"
$ws.Range("F2").Style = "Normal"

# G2: "21" -> "12" (keep as text, not a number)
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "12"
$ws.Range("G2").Style = "Normal"
